$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths for A and B (compensating for Excel's internal padding so the
# saved XML width attribute ends up as 21.5 / 18 respectively)
$ws.Columns.Item(1).ColumnWidth = 20.666666666666668
$ws.Columns.Item(2).ColumnWidth = 17.166666666666668

# Row 3
$ws.Cells.Item(3, 1).Value = "Neuhaus im Solling"
$ws.Cells.Item(3, 2).Value = "http://d-nb.info/gnd/2020253-2"
$ws.Cells.Item(3, 3).Value = "http://vocab.getty.edu/tgn/1040455"

# Row 4
$ws.Cells.Item(4, 1).Value = "Jerusalem"
$ws.Cells.Item(4, 2).Value = "http://d-nb.info/gnd/4028586-8"
$ws.Cells.Item(4, 3).Value = "http://vocab.getty.edu/tgn/7001371"

# Row 5
$ws.Cells.Item(5, 1).Value = "Lübeck"
$ws.Cells.Item(5, 2).Value = "http://d-nb.info/gnd/4036483-5"
$ws.Cells.Item(5, 3).Value = "http://vocab.getty.edu/tgn/7012327"

# Row 6
$ws.Cells.Item(6, 1).Value = "Celle"
$ws.Cells.Item(6, 2).Value = "http://d-nb.info/gnd/4009657-9"
$ws.Cells.Item(6, 3).Value = "http://vocab.getty.edu/tgn/7005317"

# Row 7 (B then A then C)
$ws.Cells.Item(7, 2).Value = "http://d-nb.info/gnd/4068038-1"
$ws.Cells.Item(7, 1).Value = "Zürich"
$ws.Cells.Item(7, 3).Value = "http://vocab.getty.edu/tgn/7007302"

# Row 8
$ws.Cells.Item(8, 1).Value = "Herne"
$ws.Cells.Item(8, 2).Value = "http://d-nb.info/gnd/4024544-5"
$ws.Cells.Item(8, 3).Value = "http://vocab.getty.edu/tgn/1039566"

# Row 9
$ws.Cells.Item(9, 1).Value = "Leipzig"
$ws.Cells.Item(9, 2).Value = "http://d-nb.info/gnd/4035206-7"
$ws.Cells.Item(9, 3).Value = "http://vocab.getty.edu/tgn/7012329"

# Row 10 (C then A then B)
$ws.Cells.Item(10, 3).Value = "http://vocab.getty.edu/tgn/7018159"
$ws.Cells.Item(10, 1).Value = "Venice"
$ws.Cells.Item(10, 2).Value = "http://d-nb.info/gnd/4062501-1:"

# Row 11 (C then A then B)
$ws.Cells.Item(11, 3).Value = "http://vocab.getty.edu/tgn/7004446"
$ws.Cells.Item(11, 1).Value = "Köln"
$ws.Cells.Item(11, 2).Value = "http://d-nb.info/gnd/4031483-2"

# Row 12
$ws.Cells.Item(12, 1).Value = "Kopenhagen"
$ws.Cells.Item(12, 2).Value = "http://d-nb.info/gnd/4032399-7"
$ws.Cells.Item(12, 3).Value = "http://vocab.getty.edu/tgn/7003474"

# Row 13
$ws.Cells.Item(13, 1).Value = "Bern"
$ws.Cells.Item(13, 2).Value = "http://d-nb.info/gnd/4005762-8"
$ws.Cells.Item(13, 3).Value = "http://vocab.getty.edu/tgn/7007557"

$ws.Range("A14").Select()
